# Update the weekly Achicoria price records (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) for rows 3-25 on
# the single data sheet, matching the refreshed "Fruta / hortaliza, semanal"
# data pull. Columns: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, P=Precio $/Kg.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value = 44371
$ws.Cells.Item(3, 10).Value = 34
$ws.Cells.Item(3, 11).Value = 5500
$ws.Cells.Item(3, 12).Value = 6000
$ws.Cells.Item(3, 13).Value = 5750
$ws.Cells.Item(3, 16).Value = 359

# Row 4 is unchanged in the source diff, so it is intentionally skipped.

# Row 5
$ws.Cells.Item(5, 4).Value = 44477
$ws.Cells.Item(5, 10).Value = 25
$ws.Cells.Item(5, 11).Value = 6000
$ws.Cells.Item(5, 12).Value = 6000
$ws.Cells.Item(5, 13).Value = 6000
$ws.Cells.Item(5, 16).Value = 375

# Row 6
$ws.Cells.Item(6, 4).Value = 44308
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(6, 11).Value = 6000
$ws.Cells.Item(6, 12).Value = 6000
$ws.Cells.Item(6, 13).Value = 6000
$ws.Cells.Item(6, 16).Value = 375

# Row 7
$ws.Cells.Item(7, 4).Value = 44376
$ws.Cells.Item(7, 10).Value = 43
$ws.Cells.Item(7, 11).Value = 4500
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = 4756
$ws.Cells.Item(7, 16).Value = 297

# Row 8
$ws.Cells.Item(8, 4).Value = 44306
$ws.Cells.Item(8, 10).Value = 50
$ws.Cells.Item(8, 11).Value = 6000
$ws.Cells.Item(8, 12).Value = 6000
$ws.Cells.Item(8, 13).Value = 6000
$ws.Cells.Item(8, 16).Value = 375

# Row 9
$ws.Cells.Item(9, 4).Value = 44467
$ws.Cells.Item(9, 10).Value = 52
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = 5500
$ws.Cells.Item(9, 16).Value = 344

# Row 10
$ws.Cells.Item(10, 4).Value = 44341
$ws.Cells.Item(10, 10).Value = 51
$ws.Cells.Item(10, 11).Value = 5500
$ws.Cells.Item(10, 12).Value = 6000
$ws.Cells.Item(10, 13).Value = 5755
$ws.Cells.Item(10, 16).Value = 360

# Row 11
$ws.Cells.Item(11, 4).Value = 44442
$ws.Cells.Item(11, 10).Value = 25
$ws.Cells.Item(11, 11).Value = 6000
$ws.Cells.Item(11, 12).Value = 7000
$ws.Cells.Item(11, 13).Value = 6480
$ws.Cells.Item(11, 16).Value = 405

# Row 12
$ws.Cells.Item(12, 4).Value = 44474
$ws.Cells.Item(12, 10).Value = 52
$ws.Cells.Item(12, 11).Value = 5000
$ws.Cells.Item(12, 12).Value = 6000
$ws.Cells.Item(12, 13).Value = 5500
$ws.Cells.Item(12, 16).Value = 344

# Row 13
$ws.Cells.Item(13, 4).Value = 44350
$ws.Cells.Item(13, 10).Value = 25
$ws.Cells.Item(13, 11).Value = 6000
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 6000
$ws.Cells.Item(13, 16).Value = 375

# Row 14
$ws.Cells.Item(14, 4).Value = 44328
$ws.Cells.Item(14, 10).Value = 160
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 6000
$ws.Cells.Item(14, 16).Value = 375

# Row 15
$ws.Cells.Item(15, 4).Value = 44589
$ws.Cells.Item(15, 10).Value = 52
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 8000
$ws.Cells.Item(15, 16).Value = 500

# Row 16
$ws.Cells.Item(16, 4).Value = 44363
$ws.Cells.Item(16, 10).Value = 160
$ws.Cells.Item(16, 11).Value = 5500
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = 5750
$ws.Cells.Item(16, 16).Value = 359

# Row 17
$ws.Cells.Item(17, 4).Value = 44582
$ws.Cells.Item(17, 10).Value = 52
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 7000
$ws.Cells.Item(17, 16).Value = 438

# Row 18
$ws.Cells.Item(18, 4).Value = 44358
$ws.Cells.Item(18, 10).Value = 52
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 6000
$ws.Cells.Item(18, 13).Value = 6000
$ws.Cells.Item(18, 16).Value = 375

# Row 19
$ws.Cells.Item(19, 4).Value = 44313
$ws.Cells.Item(19, 10).Value = 34
$ws.Cells.Item(19, 11).Value = 6000
$ws.Cells.Item(19, 12).Value = 6000
$ws.Cells.Item(19, 13).Value = 6000
$ws.Cells.Item(19, 16).Value = 375

# Row 20
$ws.Cells.Item(20, 4).Value = 44438
$ws.Cells.Item(20, 10).Value = 34
$ws.Cells.Item(20, 11).Value = 5000
$ws.Cells.Item(20, 12).Value = 6000
$ws.Cells.Item(20, 13).Value = 5500
$ws.Cells.Item(20, 16).Value = 344

# Row 21
$ws.Cells.Item(21, 4).Value = 44573
$ws.Cells.Item(21, 10).Value = 34
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 8000
$ws.Cells.Item(21, 16).Value = 500

# Row 22
$ws.Cells.Item(22, 4).Value = 44403
$ws.Cells.Item(22, 10).Value = 43
$ws.Cells.Item(22, 11).Value = 6000
$ws.Cells.Item(22, 12).Value = 6000
$ws.Cells.Item(22, 13).Value = 6000
$ws.Cells.Item(22, 16).Value = 375

# Row 23
$ws.Cells.Item(23, 4).Value = 44455
$ws.Cells.Item(23, 10).Value = 52
$ws.Cells.Item(23, 11).Value = 5000
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5500
$ws.Cells.Item(23, 16).Value = 344

# Row 24
$ws.Cells.Item(24, 4).Value = 44355
$ws.Cells.Item(24, 10).Value = 25
$ws.Cells.Item(24, 11).Value = 6000
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 6000
$ws.Cells.Item(24, 16).Value = 375

# Row 25
$ws.Cells.Item(25, 4).Value = 44330
$ws.Cells.Item(25, 10).Value = 120
$ws.Cells.Item(25, 11).Value = 6000
$ws.Cells.Item(25, 12).Value = 6000
$ws.Cells.Item(25, 13).Value = 6000
$ws.Cells.Item(25, 16).Value = 375
